$wb = $excel.ActiveWorkbook

# --- model_development sheet: insert a new blank row above the existing data row 2 ---
$wsModelDev = $wb.Worksheets.Item("model_development")
$wsModelDev.Rows.Item(2).Insert()

# --- historical_calibration sheet: fill in the previously-empty row 3 ---
$wsHist = $wb.Worksheets.Item("historical_calibration")

# Bring over the formatting (styles) used by row 4 so the new row matches
# the existing look (date style, label style, wrapped-text style) without
# creating brand-new style entries.
$wsHist.Range("A4:C4").Copy()
$wsHist.Range("A3:C3").PasteSpecial(-4122)

$wsHist.Range("A3").Value = 45012
$wsHist.Range("B3").Value = "Simple geology model"
$wsHist.Range("C3").Value = "With the local models I found that high vertical resolution led to a strong control on model drainage with the GHB extern boundary. I looked at using CVHM2 but the output is still not available, so ideally I should use the model output from the regional model to reflect the external boundary conditions in the local model with specified flux boundaries."
$wsHist.Rows.Item(3).RowHeight = 86.4

# Widen column C so it's no longer sharing a width with column D.
$wsHist.Columns.Item(3).ColumnWidth = 61.5

# --- Update selections / active sheet ---
# Select the new row on model_development first, then finish by making
# historical_calibration the active (tab-selected) sheet with C3 selected,
# matching the final UI state captured in the workbook.
$wsModelDev.Activate()
$wsModelDev.Rows.Item(2).Select()

$wsHist.Activate()
$wsHist.Range("C3").Select()
